$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: was CasesTab -> now ParticipantsTab, with refreshed participants query
$ws.Range("A2").Value = "ParticipantsTab"
$participantsQuery = @'
MATCH (s:study)<--(p:participant)
WHERE s.study_name in ["UCSF Database for the Advancement of JMML - Integration of Metadata with \"\"Omic\"\" Data"]
OPTIONAL MATCH (p)<--(samp:sample)
WITH s, p, apoc.coll.sort(collect(distinct coalesce(samp.sample_id, "Not specified in data"))) as samp
RETURN   
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
 coalesce(p.gender,'') as `Gender`,
 coalesce(apoc.text.join(samp, ','), '') as `Samples`
 ORDER By p.participant_id 
 LIMIT 100
'@
$ws.Range("B2").Value = $participantsQuery
$countsQuery = @'
MATCH (s:study)
WHERE 
    s.study_name in ["UCSF Database for the Advancement of JMML - Integration of Metadata with \"\"Omic\"\" Data"]
WITH s, {studies: count(distinct s)} AS counts
OPTIONAL MATCH (s)<--(f:file)
WITH s, apoc.map.merge(counts, {files: count(distinct f)}) AS counts
OPTIONAL MATCH (s)<--(p:participant)
WITH s, apoc.map.merge(counts, {participants: count(distinct p)}) AS counts
OPTIONAL MATCH (s)<--(:participant)<--(samp:sample)
WITH apoc.map.merge(counts, {samples: count(distinct samp)}) AS counts
RETURN
    counts.studies AS Studies,
    counts.participants AS Participants,
    counts.samples AS Samples,
    counts.files AS Files
'@
$ws.Range("C2").Value = $countsQuery

# Row 3: SamplesTab stays, refresh the dbExcel (counts) query
$ws.Range("A3").Value = "SamplesTab"
$samplesQuery = @'
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["UCSF Database for the Advancement of JMML - Integration of Metadata with \"\"Omic\"\" Data"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $countsQuery

# Row 4: FilesTab stays, refresh query + dbExcel (counts) query
$ws.Range("A4").Value = "FilesTab"
$filesQuery = @'
MATCH (s:study)
WHERE s.study_name in ["UCSF Database for the Advancement of JMML - Integration of Metadata with \"\"Omic\"\" Data"]
WITH s, "Not specified in data" as na
WITH s, na, {
    study_name: coalesce(s.study_name, na),
    phs_accession: coalesce(s.phs_accession, na)
} as output
OPTIONAL MATCH (f:file)-->(s)
WITH s, na, f, apoc.map.merge(output, {
    file_name: coalesce(f.file_name, na),
    file_type: coalesce(f.file_type, na)
}) as output
WITH f, na, output
OPTIONAL MATCH (f)-->(:sample)-->(p:participant)
WITH f, na, apoc.map.merge(output, {
    participant_id: coalesce(p.participant_id, na)
}) as output
OPTIONAL MATCH (f)-->(samp:sample)
WITH apoc.map.merge(output, {
    sample_id: coalesce(samp.sample_id, na)
}) as output
RETURN
    output.file_name as `File Name`,
    output.study_name as `Study Name`,
    output.phs_accession as `Accession`,
    output.participant_id as `Participant ID`,
    output.sample_id as `Sample ID`,
    output.file_type as `File Type`
ORDER BY `File Name` LIMIT 100
'@
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $countsQuery

# Update the saved selection to D4, matching the authored workbook state
$ws.Range("D4").Select()
